# Actualización automática 2025-10-21 16:30:09
$wb = $excel.ActiveWorkbook

# ---- Sheet "VENTAS POR GRUPO" ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M4").Value = 3969.56
$ws1.Range("M24").Value = 5544.61
$ws1.Range("K29").Value = 1710.36
$ws1.Range("M29").Value = 5863.94
$ws1.Range("M36").Value = 8092.71

# ---- Sheet "VENTA MENSUAL" ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 5589.55
$ws2.Range("F24").Value = 6331.21
$ws2.Range("F29").Value = 10389.8
$ws2.Range("F36").Value = 11374.2
$ws2.Range("F60").Value = 53442.47
# Stored OOXML col width comes back as ColumnWidth + 5/6, so back that off
# here to land on the target stored width of 13.
$ws2.Columns.Item(6).ColumnWidth = 13 - (5/6)

# ---- Sheet "CUMPLIMIENTO MENSUAL" ----
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D10").Value = 5471.21
$ws3.Range("E10").Value = -1590.13016465608
$ws3.Range("F10").Value = 1.409713335493695

$ws3.Range("D12").Value = 27850.65
$ws3.Range("E12").Value = 24812.47
$ws3.Range("F12").Value = 0.5288454235145962

$ws3.Range("D14").Value = 51347.98
$ws3.Range("E14").Value = 47668.52661190614
$ws3.Range("F14").Value = 0.5185800000120961
